$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Trabajadores" table's PK column: rut INT -> rut VARCHAR(10)
$ws.Range("D5").Value = "rut VARCHAR(10) [PK]"

# Update the "Trabajadores" table's FK column: id_departamento [FK] NOT NULL -> id_departamento INT [FK] NOT NULL
$ws.Range("D9").Value = "id_departamento INT [FK] NOT NULL"

# Widen column D a bit (no longer best-fit, fixed width ~31.1640625)
$ws.Columns.Item(4).ColumnWidth = 30.3

# Move the active selection from H18 to F18
$ws.Range("F18").Select()
